$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.751.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.384.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.12"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.84"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.385.39"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.396"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.964.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.369.63"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.838.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.01"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.558"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.74"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.520.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.39"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.43"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.96"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "165.97"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.415.53"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.96"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "28.10"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.40"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.24%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.501.18"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.49"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.79"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.71%  "
